$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Indicators")

# New, larger feature dictionary used for the "x" row (created dictionary)
$ws.Range("B2").Value = "Ones, S1_GRD_4TO49, S1_GRD_5TO59, S2_GRD_3TO39, S2_GRD_4TO49, S2_GRD_5TO59, S2_GRD_6TO7, S2_BEST_GRD, SchoolRegion_1, SchoolRegion_2, SchoolRegion_3, SchoolType_2, MotherEd_7, Campus_1"

# Re-run metrics after fixing the pyplot reset issue between MLR and KMeans
$ws.Range("B5").Value = 0.86486486486486491
$ws.Range("B6").Value = 0.84210526315789469
$ws.Range("B7").Value = 0.85333333333333339

# Column B keeps its best-fit sizing but now has to fit the longer dictionary text
$ws.Columns.Item(2).ColumnWidth = 190.5
